# Refresh the cryptocurrency market-data table (coin names, links,
# prices and 1h volume change) with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.479.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.99%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.849.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.54%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'241.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.33%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -2.30%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07540"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2975"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.10%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'24.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07686"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.18%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.909.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.69%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.003"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.77%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6856"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.93%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'83.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.000009785"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "'2.155.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'6.217"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.71%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'29.552.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'234.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'12.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'7.598"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'155.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.16%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1393"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.96%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.423"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.23%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.479"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.93%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.05839"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.82%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.259"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.25%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.105"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.37%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.018"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.09%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.881"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.170"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.44%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7184"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.54%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.588"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.64%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.797"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.93%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'1.237.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.39%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9121"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.88%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.73%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "'2.066.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.06%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.9998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'67.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.01%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +8.23%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000117"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.47%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.4028"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.82%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.122"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.10%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.704"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.84%  "
$ws.Range("E51").Style = "Normal"
